# 25A element - changed dV from 5% to 10%
# Column C on the "relays" sheet holds "Max. Slip Voltage [%]" (dV).
# Update every row where dV is currently 5% up to 10%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2..20
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 5) {
        $cell.Value = 10
    }
}

# Restore the style consistency for rows 19-20 (C19/C20), which previously
# carried a stray "applyNumberFormat + applyFill" style distinct from the
# rest of column C; align them with the common column style by copying the
# format from C18 (same column, already using the common style).
$ws.Range("C18").Copy()
$ws.Range("C19:C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to E11, matching the author's last interaction.
$ws.Range("E11").Select()
